$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,5,6,7,8,9,12,14,15)

$data = @{
    2 = @(13.52983389648876,10.75166404438965,11.98871058941985,16.86991607391245,21.46689239003918,12.43874642079204,19.00450643256271,9.560096254326819,16.48388625455781,17.95683166579544)
    3 = @(12.98457987743057,10.63848543835829,12.02740444689617,15.89584955866815,21.48037086845421,12.48223220185984,19.10711238487939,9.531828542791331,16.50710719070581,18.0198959568959)
    4 = @(12.63853139835659,10.56819884528809,12.05337888553044,15.26997757108489,21.49769717877561,12.51111672525247,19.17421343102896,9.516049571054596,16.52306604801823,18.06313823001819)
    5 = @(12.49487753958265,10.53937602458751,12.06452062060916,15.008197319934,21.50702595532251,12.52343632224611,19.20258835632179,9.510020980191065,16.52999794755721,18.08189326312125)
    6 = @(12.470870873221,10.5345796850817,12.0664043243123,14.96433081551589,21.5087117170944,12.52551512697196,19.20736223181017,9.509044317587389,16.53117488923941,18.08507588176683)
    7 = @(12.63660443102283,10.56781083570902,12.05352689231133,15.26647399323133,21.49781381793013,12.51128064947203,19.17459193242911,9.515966635788239,16.52315779789517,18.06338658206102)
    8 = @(13.3442772501674,10.71281647977039,12.00159204189911,16.5399640634477,21.469657406898,12.45328684352069,19.03903363093264,9.550025115543798,16.49154023785339,17.97763606790891)
    9 = @(14.63525995720682,10.99003270627679,11.91734810956468,19.00274580682531,21.48649106005578,12.35690340289024,18.80575939709204,9.629104581010742,16.44300183704827,17.84548890781045)
    10 = @(15.5163076650511,11.188254113164,11.86620329988556,20.67494806633232,21.54294929985474,12.29667966720441,18.65424575283672,9.694372352774106,16.41550293504116,17.77053690571969)
    11 = @(15.90106887217022,11.27702113471106,11.84527405481781,21.3917225636224,21.5781927199966,12.27158529150847,18.58964038271574,9.725544660153085,16.40475543108701,17.74128220559515)
    12 = @(16.04436699175297,11.31041445200272,11.83768500731629,21.65686569030329,21.59290803252267,12.26241401866619,18.56579789074315,9.737554921022511,16.400938132047,17.73090308096554)
    13 = @(16.01361343196833,11.3032327188839,11.83930447720553,21.60004134736742,21.58967803447036,12.26437446931125,18.57090510905902,9.734959250974095,16.40174903667108,17.73310728046133)
    14 = @(15.91290675290105,11.27977295459825,11.84464295524606,21.4136618050453,21.57937595251331,12.27082412066838,18.58766637101279,9.726528664865013,16.40443632287421,17.74041428444233)
    15 = @(15.85090547222894,11.26537386406808,11.84795674746335,21.29868154950795,21.57324376357547,12.27481789232889,18.59801418524116,9.721391297747278,16.40611522805531,17.74498114534793)
    16 = @(15.49083137111737,11.18242309032195,11.86761813028676,20.62722412089977,21.54083803550678,12.29836599762036,18.65855486653215,9.692364385798532,16.4162406945639,17.77254643466192)
    17 = @(15.26575760128384,11.13116211049912,11.88027853662919,20.20408069597325,21.52340376811867,12.31340183425502,18.69680156171637,9.674931864418049,16.42290303128233,17.79069889751346)
    18 = @(15.13479696348698,11.10154700976583,11.88778039629048,19.95656407809801,21.51427637820022,12.32226667104718,18.71920649901396,9.665045043830105,16.42690091762563,17.80159524789795)
    19 = @(15.09020080263581,11.09149790325935,11.89035815805374,19.87204792380568,21.5113407574045,12.32530534386481,18.72686220150882,9.661721775872978,16.42828304876734,17.8053627142271)
    20 = @(15.28987351306039,11.13663263975556,11.87890805013267,20.24955283636154,21.52516652477843,12.31177882107801,18.69268806192233,9.67677315682398,16.42217665133867,17.78871936976624)
    21 = @(15.94255265613,11.28666979652538,11.8430657821093,21.46857628470577,21.58236481387081,12.26892070352897,18.58272628775394,9.728999399831418,16.40364015425508,17.73824904827062)
    22 = @(16.35507660186748,11.38343225384324,11.82160175569077,22.22866616901552,21.62772483342332,12.24284243060667,18.51448689070626,9.764329424087052,16.39299713259029,17.7093392175775)
    23 = @(16.13621796146176,11.33191292393805,11.83287796595265,21.82633154458858,21.6027877996191,12.25658397507923,18.55057529188712,9.745365996647825,16.39854312977305,17.72439515728463)
    24 = @(15.27897556936211,11.13415986208371,11.87952695179665,20.22900810905287,21.52436679177534,12.31251189841696,18.6945464781932,9.675940286340785,16.42250452552485,17.7896128808058)
    25 = @(14.29732260586225,10.91591721842089,11.93825197447438,18.34778573295695,21.47419703928709,12.38111931793659,18.86537850860147,9.606427807574118,16.45469603067785,17.84548890781045)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $rowVals[$i]
    }
}

Write-Host "Updated $($data.Keys.Count) rows"
